$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers,
# preserving the original text representation (matches source inlineStr cells).
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = '@'
}

# Apply the updated coin data (price + 1h volume change).
$ws.Range('D2').Value = '29.121.58'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '1.833.21'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '239.82'
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('D6').Value = '0.6635'
$ws.Range('E6').Value = '  -4.70%  '
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.2953'
$ws.Range('E8').Value = '  -3.74%  '
$ws.Range('D9').Value = '0.07339'
$ws.Range('D10').Value = '22.71'
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('D11').Value = '0.07678'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '1.839.87'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').Value = '5.017'
$ws.Range('E13').Value = '  -2.55%  '
$ws.Range('D14').Value = '0.6739'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '86.28'
$ws.Range('E15').Value = '  -5.19%  '
$ws.Range('D16').Value = '6.102'
$ws.Range('E16').Value = '  -2.56%  '
$ws.Range('D17').Value = '29.122.90'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '0.000008222'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').Value = '228.09'
$ws.Range('E19').Value = '  -4.25%  '
$ws.Range('D20').Value = '12.48'
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '7.287'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '160.44'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').Value = '0.1416'
$ws.Range('E25').Value = '  -5.16%  '
$ws.Range('D26').Value = '8.656'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').Value = '17.99'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '1.501'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('D29').Value = '4.234'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '4.096'
$ws.Range('D31').Value = '1.195'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').Value = '0.05338'
$ws.Range('E32').Value = '  +4.87%  '
$ws.Range('D33').Value = '1.862'
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').Value = '0.7459'
$ws.Range('E34').Value = '  -3.58%  '
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '1.321.07'
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('E38').Value = '  -3.70%  '
$ws.Range('D39').Value = '2.712'
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('D40').Value = '0.9219'
$ws.Range('E40').Value = '  -2.73%  '
$ws.Range('D41').Value = '6.027'
$ws.Range('E41').Value = '  +4.47%  '
$ws.Range('D42').Value = '0.9985'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').Value = '103.28'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').Value = '1.980.73'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').Value = '0.5171'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('B46').Value = 'XinFinNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D46').Value = '0.07694'
$ws.Range('E46').Value = '  +14.76%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000121'
$ws.Range('E47').Value = '  -3.96%  '
$ws.Range('D48').Value = '1.758'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('D49').Value = '63.37'
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('D50').Value = '9.247'
$ws.Range('E50').Value = '  -5.41%  '
$ws.Range('D51').Value = '0.05923'
$ws.Range('E51').Value = '  -0.05%  '
